$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text tweaks (rich-text shared strings) - edit only the affected
# sub-run of text via Characters() so the rest of the string is untouched.
# ---------------------------------------------------------------------------

# A8: "Volume 31   Number  14" -> "...15" (issue number bump)
$ws.Range("A8").Characters(21, 2).Text = "15"

# C9: "Report Covering the Week  4/1/2024  Through  4/7/2024"
#     -> "...4/8/2024  Through  4/14/2024" (one week later)
$ws.Range("C9").Characters(27, 8).Text = "4/8/2024"
$ws.Range("C9").Characters(46, 8).Text = "4/14/2024"

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("M15").Value = 0

# ---------------------------------------------------------------------------
# Row 16 - C16/D16/E16/F16/G16/H16 flip between the "no data" text markers
# ("0" / "***.*") and real numbers, so number format must be (re)applied
# before writing values, by copying formats from rows that already carry
# the desired style.
# ---------------------------------------------------------------------------
$ws.Range("I15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1

$ws.Range("K15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = -100

$ws.Range("I15").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = 1

$ws.Range("K15").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = -100

$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 50
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = -62.5

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 10
$ws.Range("I17").Value = 34
$ws.Range("J17").Value = 24
$ws.Range("K17").Value = 41.666666666666
$ws.Range("L17").Value = 142.857142857143
$ws.Range("M17").Value = 183.333333333333
$ws.Range("N17").Value = 88.888888888888

# ---------------------------------------------------------------------------
# Row 18 - C18 flips from a real number back to the "0" text marker.
# ---------------------------------------------------------------------------
$ws.Range("C18").Formula = "'0"
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = -33.333333333333
$ws.Range("M18").Value = -82.35294117647
$ws.Range("N18").Value = -91.549295774647

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -75
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -37.931034482758
$ws.Range("I19").Value = 73
$ws.Range("J19").Value = 83
$ws.Range("K19").Value = -12.048192771084
$ws.Range("L19").Value = -9.876543209876
$ws.Range("M19").Value = 62.222222222222
$ws.Range("N19").Value = 78.048780487804

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 22
$ws.Range("K20").Value = -36.363636363636
$ws.Range("L20").Value = -56.25
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -92.708333333333

# ---------------------------------------------------------------------------
# Row 21 (TOTAL row)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 4
$ws.Range("E21").Value = -71.428571428571
$ws.Range("F21").Value = 33
$ws.Range("G21").Value = 47
$ws.Range("H21").Value = -29.787234042553
$ws.Range("I21").Value = 136
$ws.Range("J21").Value = 154
$ws.Range("K21").Value = -11.688311688311
$ws.Range("L21").Value = 3.030303030303
$ws.Range("M21").Value = 28.301886792452
$ws.Range("N21").Value = -60.117302052785

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 34
$ws.Range("H24").Value = -5.882352941176
$ws.Range("I24").Value = 117
$ws.Range("J24").Value = 151
$ws.Range("K24").Value = -22.51655629139
$ws.Range("L24").Value = 9.345794392523
$ws.Range("M24").Value = -13.970588235294

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 11.111111111111
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 86
$ws.Range("K25").Value = -31.395348837209
$ws.Range("L25").Value = 156.521739130435

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 38
$ws.Range("J26").Value = 56
$ws.Range("K26").Value = -32.142857142857
$ws.Range("L26").Value = -28.301886792452
$ws.Range("M26").Value = -34.482758620689
